$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.069.74"
$ws.Range("E2").Value = "  +3.17%  "
$ws.Range("D3").Value = "3.743.20"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'601.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").Value = "'169.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("D7").Value = "3.740.89"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'0.533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.91%  "
$ws.Range("E10").Value = "  +5.72%  "
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").Value = "'0.462"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").Value = "'38.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("D14").Value = "'0.0000245"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").Value = "4.367.78"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").Value = "3.747.70"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "69.061.70"
$ws.Range("E17").Value = "  +2.91%  "
$ws.Range("E18").Value = "  +3.55%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'17.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.66%  "
$ws.Range("D21").Value = "'498.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.34%  "
$ws.Range("D22").Value = "'9.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.76%  "
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("D24").Value = "'85.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000143"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.15%  "
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").Value = "'2.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("D27").Value = "'12.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.90%  "
$ws.Range("D28").Value = "'10.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +2.06%  "
$ws.Range("D31").Value = "'8.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.31%  "
$ws.Range("D32").Value = "'2.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.36%  "
$ws.Range("D33").Value = "'31.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "3.882.80"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("E35").Value = "  +2.26%  "
$ws.Range("D36").Value = "3.682.72"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("E39").Value = "  +3.07%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").Value = "'0.326"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("D42").Value = "'438.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").Value = "'49.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("E44").Value = "  +2.59%  "
$ws.Range("D45").Value = "'2.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("E46").Value = "  +3.20%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").Value = "'40.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D49").Value = "'142.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.03%  "
$ws.Range("D50").Value = "'0.0353"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("D51").Value = "2.761.88"
$ws.Range("E51").Value = "  -0.80%  "